# Updated cryptos list on Mon Feb 26 10:51:25 UTC 2024 with GitHub Actions
# Applies the per-row Price (D) and Volume(1h) (E) updates to the crypto tracker sheet.
# Values are assigned with a leading apostrophe so Excel preserves them as literal
# text (matching the original inlineStr cells) instead of re-interpreting numeric-
# looking strings (e.g. "50.00", "3.064.52") as numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.197.16"
$ws.Range("E2").Value = "'  -0.75%  "
$ws.Range("D3").Value = "'3.063.35"
$ws.Range("E3").Value = "'  +1.30%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'388.07"
$ws.Range("E5").Value = "'  +2.23%  "
$ws.Range("D6").Value = "'102.16"
$ws.Range("E6").Value = "'  -0.11%  "
$ws.Range("E7").Value = "'  -1.82%  "
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E9").Value = "'  -1.64%  "
$ws.Range("D10").Value = "'36.76"
$ws.Range("E10").Value = "'  +0.37%  "
$ws.Range("E11").Value = "'  +0.44%  "
$ws.Range("D12").Value = "'0.0849"
$ws.Range("E12").Value = "'  -1.30%  "
$ws.Range("D13").Value = "'3.548.36"
$ws.Range("E13").Value = "'  +1.32%  "
$ws.Range("E14").Value = "'  -1.21%  "
$ws.Range("D15").Value = "'7.67"
$ws.Range("E15").Value = "'  -0.67%  "
$ws.Range("D16").Value = "'3.060.02"
$ws.Range("E16").Value = "'  +1.11%  "
$ws.Range("D17").Value = "'0.994"
$ws.Range("E17").Value = "'  +2.26%  "
$ws.Range("D18").Value = "'10.59"
$ws.Range("E18").Value = "'  +0.20%  "
$ws.Range("D19").Value = "'51.217.46"
$ws.Range("E20").Value = "'  +2.77%  "
$ws.Range("E21").Value = "'  -1.29%  "
$ws.Range("E22").Value = "'  -0.61%  "
$ws.Range("D23").Value = "'69.67"
$ws.Range("D24").Value = "'264.17"
$ws.Range("E24").Value = "'  -1.18%  "
$ws.Range("D25").Value = "'3.12"
$ws.Range("E25").Value = "'  -0.99%  "
$ws.Range("E26").Value = "'  -6.95%  "
$ws.Range("D27").Value = "'26.84"
$ws.Range("E27").Value = "'  +2.65%  "
$ws.Range("D28").Value = "'7.22"
$ws.Range("E28").Value = "'  -4.03%  "
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("E30").Value = "'  -5.02%  "
$ws.Range("D31").Value = "'0.104"
$ws.Range("E31").Value = "'  -3.34%  "
$ws.Range("D32").Value = "'10.43"
$ws.Range("E32").Value = "'  +1.64%  "
$ws.Range("D33").Value = "'35.42"
$ws.Range("E33").Value = "'  +4.21%  "
$ws.Range("D34").Value = "'0.0472"
$ws.Range("E34").Value = "'  +5.12%  "
$ws.Range("E35").Value = "'  +2.35%  "
$ws.Range("D36").Value = "'50.00"
$ws.Range("E36").Value = "'  -1.16%  "
$ws.Range("E37").Value = "'  -0.04%  "
$ws.Range("E38").Value = "'  +1.38%  "
$ws.Range("E39").Value = "'  -1.28%  "
$ws.Range("D40").Value = "'130.81"
$ws.Range("E40").Value = "'  +1.35%  "
$ws.Range("D41").Value = "'16.47"
$ws.Range("E41").Value = "'  -3.45%  "
$ws.Range("E42").Value = "'  -1.61%  "
$ws.Range("E43").Value = "'  -0.85%  "
$ws.Range("E44").Value = "'  -1.86%  "
$ws.Range("D45").Value = "'3.72"
$ws.Range("E45").Value = "'  -0.33%  "
$ws.Range("D46").Value = "'21.68"
$ws.Range("E46").Value = "'  +0.62%  "
$ws.Range("D47").Value = "'2.51"
$ws.Range("E47").Value = "'  +3.17%  "
$ws.Range("E48").Value = "'  -0.36%  "
$ws.Range("D49").Value = "'2.068.24"
$ws.Range("E49").Value = "'  +2.16%  "
$ws.Range("D50").Value = "'0.0325"
$ws.Range("E50").Value = "'  +3.79%  "
$ws.Range("D51").Value = "'0.901"
$ws.Range("E51").Value = "'  +13.76%  "
